# Weekly update: add two new records (rows) at the top of the data block
# (row 16/17), pushing all subsequent rows down by two positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 16 (existing rows 16-90 shift to 18-92)
$ws.Rows("16:17").Insert()

# ---- New row 16 ----
$ws.Range("A16").Value2 = 11
$ws.Range("B16").Value2 = "Vega Monumental Concepción"
$ws.Range("C16").Value2 = "Bíobío"
$ws.Range("D16").Value2 = 44565
$ws.Range("D16").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E16").Value2 = 8
$ws.Range("F16").Value2 = "Fruta"
$ws.Range("G16").Value2 = 100103
$ws.Range("H16").Value2 = "Frutos de hueso (carozo)"
$ws.Range("I16").Value2 = 100103001
$ws.Range("J16").Value2 = "Cereza"
$ws.Range("K16").Value2 = "Lapins"
$ws.Range("L16").Value2 = "Primera"
$ws.Range("M16").Value2 = 250
$ws.Range("N16").Value2 = 4000
$ws.Range("O16").Value2 = 5000
$ws.Range("P16").Value2 = 4400
$ws.Range("Q16").Value2 = "`$/caja 10 kilos"
$ws.Range("R16").Value2 = "Provincia de Curicó"
$ws.Range("S16").Value2 = 440
$ws.Range("T16").Value2 = 10

# ---- New row 17 ----
$ws.Range("A17").Value2 = 11
$ws.Range("B17").Value2 = "Vega Monumental Concepción"
$ws.Range("C17").Value2 = "Bíobío"
$ws.Range("D17").Value2 = 44565
$ws.Range("D17").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E17").Value2 = 8
$ws.Range("F17").Value2 = "Fruta"
$ws.Range("G17").Value2 = 100103
$ws.Range("H17").Value2 = "Frutos de hueso (carozo)"
$ws.Range("I17").Value2 = 100103001
$ws.Range("J17").Value2 = "Cereza"
$ws.Range("K17").Value2 = "Rainier"
$ws.Range("L17").Value2 = "Primera"
$ws.Range("M17").Value2 = 220
$ws.Range("N17").Value2 = 5000
$ws.Range("O17").Value2 = 5500
$ws.Range("P17").Value2 = 5227
$ws.Range("Q17").Value2 = "`$/bandeja 10 kilos"
$ws.Range("R17").Value2 = "Provincia de Curicó"
$ws.Range("S17").Value2 = 523
$ws.Range("T17").Value2 = 10
